# Apply the Mon Jul  1 08:50:38 UTC 2024 cryptos-list refresh.
# Each Range.Value assignment is prefixed with a literal apostrophe so
# Excel's type inference treats numeric-looking strings (e.g. "1.00",
# "3.476.69") as literal text instead of coercing them to a Number and
# silently dropping formatting like trailing zeros. Style is then reset
# to "Normal" so the one-time quote-prefix flag does not linger on the
# cell format (keeps cell styling identical to the un-touched cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.879.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.36%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.476.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.47%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'583.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.16%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'147.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.33%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.90%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'7.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.59%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.92%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +2.38%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'4.081.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.72%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'29.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.19%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("E14").Value = "'  -0.08%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.491.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.01%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.54%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'63.132.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.71%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("E18").Value = "'  +2.53%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  +4.74%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'9.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.46%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'387.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.11%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  +1.46%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'74.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.50%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  -0.02%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'3.630.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.84%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +2.77%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.181"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -5.90%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'7.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.02%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.17%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.60%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +0.04%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'23.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.46%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  +4.91%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +2.17%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  +21.39%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'171.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.00%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'  +6.59%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'3.521.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.82%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.0768"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.30%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.35%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'OKB"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'42.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.43%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'4.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.62%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "'Stacks"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.11%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'ONDO"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.36%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'2.619.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +6.23%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'23.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.18%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +9.11%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  +1.01%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'FirstDigitalUSD"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.08%  "
$ws.Range("E51").Style = "Normal"
